$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Characters(21, 1).Text = "10"
$ws.Range("C9").Characters(27, 9).Text = "3/3/2025"
$ws.Range("C9").Characters(47, 8).Text = "3/9/2025"

# --- Column E width will auto-adjust; set explicitly to match target ---

# --- Fix cell styles before/while writing new values where type changes ---
# Row 16: C16 text -> number
$ws.Range("D16").Copy($ws.Range("C16"))

# Row 22: C22, D22 text -> number; E22 text -> %chg number
$ws.Range("F22").Copy($ws.Range("C22"))
$ws.Range("F22").Copy($ws.Range("D22"))
$ws.Range("H22").Copy($ws.Range("E22"))

# Row 31: C31, D31 number -> text "0"; E31 number -> text "***.*"
$ws.Range("C23").Copy($ws.Range("C31"))
$ws.Range("D23").Copy($ws.Range("D31"))
$ws.Range("E23").Copy($ws.Range("E31"))

# --- Row value updates ---
# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 17
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = 21.428571428571
$ws.Range("L16").Value = 54.545454545454
$ws.Range("M16").Value = -39.285714285714
$ws.Range("N16").Value = -88.590604026845

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 160
$ws.Range("I17").Value = 31
$ws.Range("J17").Value = 15
$ws.Range("K17").Value = 106.666666666667
$ws.Range("L17").Value = 158.333333333333
$ws.Range("M17").Value = 210
$ws.Range("N17").Value = 106.666666666667

# Row 18
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 28.571428571428
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = 111.111111111111
$ws.Range("L18").Value = 18.75
$ws.Range("M18").Value = 65.217391304347
$ws.Range("N18").Value = -86.851211072664

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -30
$ws.Range("I19").Value = 63
$ws.Range("J19").Value = 69
$ws.Range("K19").Value = -8.695652173913
$ws.Range("L19").Value = -35.051546391752
$ws.Range("M19").Value = -13.698630136986
$ws.Range("N19").Value = -64.804469273743

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -46.666666666666
$ws.Range("I20").Value = 24
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = -4
$ws.Range("L20").Value = -7.692307692307
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -96.433878157503

# Row 21
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 64
$ws.Range("H21").Value = -14.0625
$ws.Range("I21").Value = 175
$ws.Range("J21").Value = 142
$ws.Range("K21").Value = 23.239436619718
$ws.Range("L21").Value = -2.234636871508
$ws.Range("M21").Value = 10.759493670886
$ws.Range("N21").Value = -86.610558530987

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = 33.333333333333

# Row 24
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 53.846153846153
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 153
$ws.Range("H24").Value = 22.222222222222
$ws.Range("I24").Value = 406
$ws.Range("J24").Value = 323
$ws.Range("K24").Value = 25.696594427244
$ws.Range("L24").Value = 35.333333333333
$ws.Range("M24").Value = 153.75

# Row 25
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 43.478260869565
$ws.Range("F25").Value = 155
$ws.Range("H25").Value = 20.155038759689
$ws.Range("I25").Value = 335
$ws.Range("J25").Value = 257
$ws.Range("K25").Value = 30.350194552529
$ws.Range("L25").Value = 60.287081339712

# Row 26
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 3.571428571428
$ws.Range("I26").Value = 66
$ws.Range("J26").Value = 47
$ws.Range("K26").Value = 40.425531914893
$ws.Range("L26").Value = 83.333333333333
$ws.Range("M26").Value = 46.666666666666

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 11
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = 83.333333333333
$ws.Range("L28").Value = 120
